# Adiciona colunas "modelo" e "politica" antes da coluna "full",
# preenchendo os valores correspondentes e atualizando os links
# (novo tracking_id) na planilha "acessorios web".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere duas novas colunas nas posições C e D (a antiga coluna "full"
# vira "E", "tipo" vira "F" e "link" vira "G").
$ws.Range("C1:D1").EntireColumn.Insert()

# Cabecalhos das novas colunas
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Dados da coluna "modelo" (C) e "politica" (D) por linha
$modelo = @{
    2 = "Modelo identificado mas fora do range de preco"
    3 = "Modelo identificado mas fora do range de preco"
    4 = "Sem Modelo"
    5 = "FONTE 70A"
    6 = "Sem Modelo"
    7 = "FONTE 70A"
    8 = "Sem Modelo"
    9 = "Sem Modelo"
}

$politica = @{
    2 = ""
    3 = ""
    4 = ""
    5 = "Igual"
    6 = ""
    7 = "Igual"
    8 = ""
    9 = ""
}

foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = $modelo[$row]
    $ws.Cells.Item($row, 4).Value = $politica[$row]
}

# Atualiza o tracking_id dos links (coluna G apos a insercao das colunas)
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 7)
    $linkValue = $cell.Value()
    $cell.Value = $linkValue.Replace("edfb7108-c6d3-4e7b-a86a-c080a604d1b1", "4e91c180-f7d0-4354-88fc-2d08e064eee4")
}

# Normaliza o texto da coluna "tipo" (coluna F apos a insercao) para
# minusculas, conforme o novo padrao.
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 6)
    $tipoValue = $cell.Value()
    $cell.Value = $tipoValue.ToLower()
}
